$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet (as reflected by the workbook's sheet name / tab)
$ws.Name = "2023_07_05 16_46"

# Update data cells per diff

# Row 2
$ws.Range("D2").Value = -72
$ws.Range("I2").Value = 24597
$ws.Range("J2").Value = 40347
$ws.Range("K2").Value = 24509

# Row 3
$ws.Range("D3").Value = -56
$ws.Range("I3").Value = 23773
$ws.Range("J3").Value = 38686
$ws.Range("K3").Value = 23655

# Row 4
$ws.Range("D4").Value = -56
$ws.Range("G4").Value = 26
$ws.Range("I4").Value = 24526
$ws.Range("J4").Value = 40215
$ws.Range("K4").Value = 24420

# Row 5
$ws.Range("D5").Value = -56
$ws.Range("I5").Value = 24571
$ws.Range("J5").Value = 40303
$ws.Range("K5").Value = 24471

# Row 6
$ws.Range("D6").Value = -68
$ws.Range("G6").Value = 25
$ws.Range("I6").Value = 24472
$ws.Range("J6").Value = 40101
$ws.Range("K6").Value = 24372

# Row 7
$ws.Range("D7").Value = -64
$ws.Range("G7").Value = 23
$ws.Range("I7").Value = 24615
$ws.Range("J7").Value = 40384
$ws.Range("K7").Value = 24527

# Row 8
$ws.Range("D8").Value = -58
$ws.Range("G8").Value = 24
$ws.Range("I8").Value = 24540
$ws.Range("J8").Value = 40235

# Row 9
$ws.Range("D9").Value = -64
$ws.Range("I9").Value = 24205
$ws.Range("J9").Value = 39552
$ws.Range("K9").Value = 24111

# Row 10
$ws.Range("G10").Value = 26
$ws.Range("I10").Value = 23862
$ws.Range("J10").Value = 38860
$ws.Range("K10").Value = 23756

# Row 11
$ws.Range("D11").Value = -64
$ws.Range("G11").Value = 25
$ws.Range("I11").Value = 24310
$ws.Range("J11").Value = 39770

$wb.Save()
